$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header label for the inflation driver column (E)
$ws.Range("E1").Value = "EURINFL"

# New row label for the inflation driver (row 5)
$ws.Range("A5").Value = "EURINFL"

# Correlation values for the new EURINFL row/column
$ws.Range("B5").Value = -0.2
$ws.Range("C5").Value = -0.3
$ws.Range("D5").Value = 0.03
$ws.Range("E5").Value = 1

# Mirror formulas in column E (symmetric correlation matrix)
$ws.Range("E2").Formula = "=B5"
$ws.Range("E3").Formula = "=C5"
$ws.Range("E4").Formula = "=D5"

# Update selection to match the target state
$ws.Range("E5").Select()
